$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$v1 = $ws.Range("I8").Value
$v2 = $ws.Range("I10").Value
$v3 = $ws.Range("L14").Value
$v4 = $ws.Range("N20").Formula
Write-Host "I8=$v1"
Write-Host "I10=$v2"
Write-Host "L14=$v3"
Write-Host "N20=$v4"
